$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - reorder category labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "living_rooms_2"

# Data rows 2-7 - one-hot encoded rows, re-permuted
$data = @(
    @(0,0,0,0,1,0),
    @(0,0,0,1,0,0),
    @(0,1,0,0,0,0),
    @(0,0,1,0,0,0),
    @(1,0,0,0,0,0),
    @(0,0,0,0,0,1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $data[$i][$j]
    }
}
